$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

# --- Phase 1: write the final text for every populated row in column A -------
# A code "PYGVMDQPPHKYLLYD " was inserted at row 3 (pushing the rest down by
# one), and four new codes were appended after the previous last code.
$ws.Cells.Item(1,1).Value  = "codes for  only Gc"
$ws.Cells.Item(2,1).Value  = "RDYHRJYOZFRIKDDT "
$ws.Cells.Item(3,1).Value  = "PYGVMDQPPHKYLLYD "
$ws.Cells.Item(4,1).Value  = "TSYTLQMVLSJJYHRP "
$ws.Cells.Item(5,1).Value  = "OTLSMRSOFPWDKWKQ "
$ws.Cells.Item(6,1).Value  = "OTLSMRSOFPWDKWKQ "
$ws.Cells.Item(7,1).Value  = "OTLSMRSOFPWDKWKQ "
$ws.Cells.Item(8,1).Value  = "QPVFHTAVILQYJRJD "
$ws.Cells.Item(9,1).Value  = "KDHVRQCPJCFMFFAY "
$ws.Cells.Item(10,1).Value = "SPQKTWFWSTPIOQOO "
$ws.Cells.Item(11,1).Value = "ZLQYYISHRYODIRQG "
$ws.Cells.Item(12,1).Value = "WTCIQJTMWWJSCSFH "
$ws.Cells.Item(13,1).Value = "VYKZILJIYJSYKHPR "
$ws.Cells.Item(14,1).Value = "JJTYTAFLPKCRHRYO "

# --- Phase 2: formatting ------------------------------------------------------
# Two bold-font flavours are used in this sheet: a "plain" bold (no explicit
# colour) and a "dark" bold (colour #1C1C1C). Re-apply the matching one to
# every cell so the style survives the row re-shuffle below intact. Source
# cells with a known-good flavour are copied via PasteSpecial(Formats) so the
# workbook's existing font/xf entries are reused instead of minting new ones.

# "Plain" bold source: header cell A1 already has this exact style.
$ws.Cells.Item(1,1).Copy()
$plainTargets = @(2,3,4,5,6,7,9,10,11,12,14,15,16,17,18)
foreach ($r in $plainTargets) {
    $ws.Cells.Item($r,1).PasteSpecial($xlPasteFormats)
}

# "Dark" bold source: A12 already has this exact style (blank cell before
# this edit, s="5" in the original workbook).
$ws.Cells.Item(12,1).Copy()
$darkTargets = @(8,13)
foreach ($r in $darkTargets) {
    $ws.Cells.Item($r,1).PasteSpecial($xlPasteFormats)
}

# Re-apply the "plain" flavour to A12 itself last, since it was (temporarily)
# used above only as a format source and must end up plain like the rest of
# the newly appended codes.
$ws.Cells.Item(1,1).Copy()
$ws.Cells.Item(12,1).PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# Move the active selection to A14, matching the edited workbook.
$ws.Range("A14").Select() | Out-Null

Write-Output "done"
